$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from SCD0226 to SCD0015
$ws.Name = "SCD0015"

# Update the TC_ID cell (B2) value
$ws.Range("B2").Value = "SCD0015-004"

# Move the selection to B3 with A2 as the top-left visible cell, matching
# the saved view state in the edited workbook.
$ws.Range("A2").Select
$ws.Range("B3").Select
